# Library_hbrown_08.22.19.xlsx — "updated files to conform to standard"
#
# Column G (index1Name) held plain numbers 1..26 for rows 2..27. These are
# replaced with standardized text labels "Index1_1" .. "Index1_26".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 26; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 7).Value = "Index1_$i"
}

# Selection moved to the column that was just edited.
$ws.Range("G2:G27").Select()
